$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-02-15 Thursday" "2024-02-16 Friday"

Replace-Text "48×92=" "19×52="
Replace-Text "96×41=" "27×27="
Replace-Text "50×80=" "17×56="
Replace-Text "93×30=" "44×60="
Replace-Text "34×14=" "63×28="
Replace-Text "57×45=" "91×54="
Replace-Text "22×72=" "85×32="
Replace-Text "19×66=" "65×65="
Replace-Text "70×42=" "64×20="
Replace-Text "73×26=" "32×18="
Replace-Text "92×34=" "56×43="
Replace-Text "69×22=" "69×12="
Replace-Text "32×64=" "39×31="
Replace-Text "89×31=" "13×81="
Replace-Text "58×93=" "39×50="
Replace-Text "52×92=" "81×58="
Replace-Text "50×33=" "36×38="
Replace-Text "95×18=" "15×60="
Replace-Text "54×63=" "39×22="
Replace-Text "18×56=" "37×57="
Replace-Text "37×43=" "98×82="
Replace-Text "99×24=" "65×73="
Replace-Text "28×97=" "12×81="
Replace-Text "18×84=" "13×26="
Replace-Text "50×83=" "82×54="
